$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start clean: wipe all existing cell content/formatting on the sheet.
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Create a throwaway named cell style that only carries the Arial-9
# font (no explicit number format). Applying it, then deleting the
# named style, leaves the underlying cell format (cellXfs) entry
# behind as "applyFont" only (no applyNumberFormat) -- matching the
# workbook's style table layout after the edit.
# ---------------------------------------------------------------------
$hdrStyle = $wb.Styles.Add("__HdrTmp__")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9

# ---------------------------------------------------------------------
# Row 1: column headers
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1:K1").Style = "__HdrTmp__"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Drop the temporary named style now that its format has been stamped
# onto F1:K1 -- keeps cellStyles/cellStyleXfs back at their original
# counts while the new cellXfs entry survives.
$wb.Styles.Item("__HdrTmp__").Delete()

# ---------------------------------------------------------------------
# Data rows 2-7. Each plant occupies one row with (possibly) some
# trailing columns left blank -- only format/write the cells that
# actually hold a value so no stray empty-but-styled cells appear.
# ---------------------------------------------------------------------
$plants = @(
    @{ idx=1; idx2=108900; name="Wunderklingen"; ds=1895; de=1968; m3s=5.5;   mw1=0.42;  mw2=0.41;  gwhw=1.4;   gwhs=1;     gwhy=2.4   },
    @{ idx=2; idx2=106300; name="Engeweiher";    ds=1909; de=1993; m3s=4;     mw1=5;     mw2=5;     gwhw=$null; gwhs=$null; gwhy=$null },
    @{ idx=3; idx2=108700; name="Eglisau";       ds=1920; de=1927; m3s=400;   mw1=10.82; mw2=10.34; gwhw=37.49; gwhs=39.59; gwhy=77.08 },
    @{ idx=4; idx2=106400; name="Neuhausen";     ds=1951; de=$null; m3s=25;   mw1=2.32;  mw2=2.2;   gwhw=9.95;  gwhs=9.95;  gwhy=19.9  },
    @{ idx=5; idx2=106500; name="Rheinau";       ds=1956; de=2005; m3s=400;   mw1=2.98;  mw2=2.92;  gwhw=6.39;  gwhs=13.19; gwhy=19.58 },
    @{ idx=6; idx2=106200; name="Schaffhausen";  ds=1964; de=$null; m3s=500;  mw1=22.57; mw2=19.84; gwhw=62.06; gwhs=73.64; gwhy=135.7 }
)

$row = 2
foreach ($p in $plants) {
    $ws.Cells.Item($row, 1).Font.Size = 9
    $ws.Cells.Item($row, 1).NumberFormat = "0"
    $ws.Cells.Item($row, 1).Value = $p.idx

    $ws.Cells.Item($row, 2).Font.Size = 9
    $ws.Cells.Item($row, 2).NumberFormat = "0"
    $ws.Cells.Item($row, 2).Value = $p.idx2

    $ws.Cells.Item($row, 3).Font.Size = 9
    $ws.Cells.Item($row, 3).Value = $p.name

    $ws.Cells.Item($row, 4).Font.Size = 9
    $ws.Cells.Item($row, 4).NumberFormat = "0"
    $ws.Cells.Item($row, 4).Value = $p.ds

    if ($p.de -ne $null) {
        $ws.Cells.Item($row, 5).Font.Size = 9
        $ws.Cells.Item($row, 5).NumberFormat = "0"
        $ws.Cells.Item($row, 5).Value = $p.de
    }

    $ws.Cells.Item($row, 6).Font.Size = 9
    $ws.Cells.Item($row, 6).NumberFormat = "0.00"
    $ws.Cells.Item($row, 6).Value = $p.m3s

    $ws.Cells.Item($row, 7).Font.Size = 9
    $ws.Cells.Item($row, 7).NumberFormat = "0.00"
    $ws.Cells.Item($row, 7).Value = $p.mw1

    $ws.Cells.Item($row, 8).Font.Size = 9
    $ws.Cells.Item($row, 8).NumberFormat = "0.00"
    $ws.Cells.Item($row, 8).Value = $p.mw2

    if ($p.gwhw -ne $null) {
        $ws.Cells.Item($row, 9).Font.Size = 9
        $ws.Cells.Item($row, 9).NumberFormat = "0.00"
        $ws.Cells.Item($row, 9).Value = $p.gwhw

        $ws.Cells.Item($row, 10).Font.Size = 9
        $ws.Cells.Item($row, 10).NumberFormat = "0.00"
        $ws.Cells.Item($row, 10).Value = $p.gwhs

        $ws.Cells.Item($row, 11).Font.Size = 9
        $ws.Cells.Item($row, 11).NumberFormat = "0.00"
        $ws.Cells.Item($row, 11).Value = $p.gwhy
    }

    $row = $row + 1
}

# Match the saved selection state from the authored workbook.
$ws.Range("A4:K4").Select()
